# daily auto push: 2026-01-22 09:44 UTC
# Insert one new timestamped reading row into the "sei1" daily log sheet.
# The new row (2026/01/22, 木, 16, 182) belongs in the existing 2026/01/22
# block, pushing every following row down by one (old row 704 -> 705, ...,
# old row 745 -> 746) and growing the used range from A1:D745 to A1:D746.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room: push rows 704..745 down to 705..746.
$ws.Rows.Item(704).Insert()

# Populate the newly inserted row 704.
# The date column stores plain text (e.g. "2026/01/22"), not a real date
# serial, so force text entry with a leading apostrophe and then strip the
# resulting "quote prefix" formatting back to the sheet's default style so
# the cell matches its style-less neighbours.
$ws.Range("A704").Value = "'2026/01/22"
$ws.Range("A704").Style = "Normal"
$ws.Range("B704").Value = "木"
$ws.Range("C704").Value = 16
$ws.Range("D704").Value = 182
